$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.051.99"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").Value = "3.442.67"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'411.71"
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("D6").Value = "'129.81"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "'0.638"
$ws.Range("E7").Value = "  +1.38%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'0.739"
$ws.Range("E9").Value = "  -2.86%  "
$ws.Range("E10").Value = "  +0.05%  "
$ws.Range("D11").Value = "'43.72"
$ws.Range("E11").Value = "  +0.28%  "
$ws.Range("E12").Value = "  +13.21%  "
$ws.Range("D13").Value = "'9.36"
$ws.Range("E13").Value = "  +4.05%  "
$ws.Range("D14").Value = "3.990.10"
$ws.Range("E14").Value = "  +0.21%  "
$ws.Range("E15").Value = "  +0.15%  "
$ws.Range("D16").Value = "'21.22"
$ws.Range("E16").Value = "  +2.33%  "
$ws.Range("D17").Value = "3.442.11"
$ws.Range("E17").Value = "  -0.20%  "
$ws.Range("D18").Value = "'12.69"
$ws.Range("E18").Value = "  +1.20%  "
$ws.Range("E19").Value = "  +1.62%  "
$ws.Range("D20").Value = "62.118.93"
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("D21").Value = "'499.72"
$ws.Range("E21").Value = "  +23.36%  "
$ws.Range("D22").Value = "'93.08"
$ws.Range("E22").Value = "  +3.20%  "
$ws.Range("D23").Value = "'3.31"
$ws.Range("E23").Value = "  +3.28%  "
$ws.Range("D24").Value = "'13.59"
$ws.Range("E24").Value = "  +1.29%  "
$ws.Range("E25").Value = "  +5.07%  "
$ws.Range("D26").Value = "'35.02"
$ws.Range("E26").Value = "  +3.64%  "
$ws.Range("D27").Value = "'9.20"
$ws.Range("E27").Value = "  +5.25%  "
$ws.Range("D28").Value = "'4.82"
$ws.Range("E28").Value = "  +0.26%  "
$ws.Range("D29").Value = "'7.64"
$ws.Range("E29").Value = "  -0.79%  "
$ws.Range("D30").Value = "'12.17"
$ws.Range("E30").Value = "  +1.75%  "
$ws.Range("E31").Value = "  -1.60%  "
$ws.Range("E32").Value = "  -2.13%  "
$ws.Range("D33").Value = "'0.169"
$ws.Range("E33").Value = "  -2.05%  "
$ws.Range("D34").Value = "'42.10"
$ws.Range("E34").Value = "  -4.48%  "
$ws.Range("D35").Value = "'59.63"
$ws.Range("E35").Value = "  +13.48%  "
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("D37").Value = "'0.0500"
$ws.Range("E37").Value = "  -0.83%  "
$ws.Range("B38").Value = "FirstDigitalUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = "  +0.14%  "
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").Value = "'3.49"
$ws.Range("E39").Value = "  +2.03%  "
$ws.Range("D40").Value = "'0.138"
$ws.Range("E40").Value = "  +4.53%  "
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").Value = "'149.21"
$ws.Range("E41").Value = "  +5.95%  "
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").Value = "'2.74"
$ws.Range("E42").Value = "  +13.20%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "'2.14"
$ws.Range("E43").Value = "  +7.67%  "
$ws.Range("D44").Value = "'2.95"
$ws.Range("E44").Value = "  +1.32%  "
$ws.Range("D45").Value = "'0.319"
$ws.Range("E45").Value = "  +1.33%  "
$ws.Range("D46").Value = "'4.33"
$ws.Range("E46").Value = "  +6.29%  "
$ws.Range("B47").Value = "ThetaToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D47").Value = "'2.36"
$ws.Range("E47").Value = "  +22.09%  "
$ws.Range("B48").Value = "Celestia"
$ws.Range("C48").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D48").Value = "'16.64"
$ws.Range("E48").Value = "  -1.31%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'23.09"
$ws.Range("E49").Value = "  +2.98%  "
$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D50").Value = "'120.29"
$ws.Range("E50").Value = "  +23.74%  "
$ws.Range("E51").Value = "  +19.48%  "
